$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -4
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 0
$ws.Range("F7").Value = 3
$ws.Range("F9").Value = -7
